$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be written as TEXT even when it looks numeric
# (e.g. "0.780" or "19.25"), without altering the cell NumberFormat/Style -
# build it via a formula that yields a text result, then convert the formula
# cell to a plain value in place (Copy + PasteSpecial values-only).
function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $escaped = $value -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$ws.Range("D2").Value = '37.472.87'
$ws.Range("E2").Value = '  +3.65%  '

$ws.Range("D3").Value = '2.067.59'
$ws.Range("E3").Value = '  +5.43%  '

$ws.Range("E4").Value = '  -0.01%  '

Set-TextValue "D5" '238.55'
$ws.Range("E5").Value = '  +4.47%  '

Set-TextValue "D6" '0.619'
$ws.Range("E6").Value = '  +4.22%  '

Set-TextValue "D7" '58.44'
$ws.Range("E7").Value = '  +9.29%  '

$ws.Range("E8").Value = '  +0.01%  '

Set-TextValue "D9" '0.383'
$ws.Range("E9").Value = '  +5.04%  '

Set-TextValue "D10" '58.03'
$ws.Range("E10").Value = '  +1.33%  '

Set-TextValue "D11" '0.0764'
$ws.Range("E11").Value = '  +3.11%  '

$ws.Range("D13").Value = '2.372.72'
$ws.Range("E13").Value = '  +5.54%  '

Set-TextValue "D14" '14.36'
$ws.Range("E14").Value = '  +4.44%  '

Set-TextValue "D15" '21.11'
$ws.Range("E15").Value = '  +6.68%  '

Set-TextValue "D16" '0.780'
$ws.Range("E16").Value = '  +4.65%  '

Set-TextValue "D17" '5.22'
$ws.Range("E17").Value = '  +5.03%  '

$ws.Range("D18").Value = '2.047.30'
$ws.Range("E18").Value = '  +4.53%  '

$ws.Range("D19").Value = '37.671.35'
$ws.Range("E19").Value = '  +4.43%  '

Set-TextValue "D20" '6.22'
$ws.Range("E20").Value = '  +24.52%  '

Set-TextValue "D21" '69.08'
$ws.Range("E21").Value = '  +2.70%  '

$ws.Range("D22").Value = '0.0₃0814'
$ws.Range("E22").Value = '  +2.58%  '

Set-TextValue "D23" '225.21'
$ws.Range("E23").Value = '  +2.00%  '

$ws.Range("E24").Value = '  -0.02%  '

Set-TextValue "D25" '2.46'
$ws.Range("E25").Value = '  +7.23%  '

Set-TextValue "D26" '2.41'
$ws.Range("E26").Value = '  +3.13%  '

Set-TextValue "D27" '163.96'
$ws.Range("E27").Value = '  +3.08%  '

Set-TextValue "D28" '8.90'
$ws.Range("E28").Value = '  +5.35%  '

Set-TextValue "D29" '1.42'
$ws.Range("E29").Value = '  +7.99%  '

Set-TextValue "D30" '19.25'
$ws.Range("E30").Value = '  +3.50%  '

Set-TextValue "D31" '0.127'
$ws.Range("E31").Value = '  +5.89%  '

$ws.Range("E32").Value = '  +3.73%  '

Set-TextValue "D33" '4.53'
$ws.Range("E33").Value = '  +5.16%  '

Set-TextValue "D34" '0.0631'
$ws.Range("E34").Value = '  +5.83%  '

Set-TextValue "D35" '2.61'
$ws.Range("E35").Value = '  +14.70%  '

Set-TextValue "D36" '4.47'
$ws.Range("E36").Value = '  +6.73%  '

$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D38" '3.36'
$ws.Range("E38").Value = '  +5.63%  '

$ws.Range("B39").Value = 'WEMIXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D39" '1.78'
$ws.Range("E39").Value = '  +0.56%  '

Set-TextValue "D40" '5.88'
$ws.Range("E40").Value = '  +14.14%  '

Set-TextValue "D41" '4.58'
$ws.Range("E41").Value = '  +30.65%  '

$ws.Range("B42").Value = 'Cronos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D42" '0.0986'
$ws.Range("E42").Value = '  +13.15%  '

$ws.Range("B43").Value = 'HuobiToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D43" '2.97'
$ws.Range("E43").Value = '  -0.21%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D44" '98.45'
$ws.Range("E44").Value = '  +13.78%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.477.73'
$ws.Range("E45").Value = '  +4.86%  '

$ws.Range("E46").Value = '  +6.65%  '

$ws.Range("E47").Value = '  +7.17%  '

Set-TextValue "D48" '15.99'
$ws.Range("E48").Value = '  +10.04%  '

Set-TextValue "D49" '1.03'
$ws.Range("E49").Value = '  +5.25%  '

Set-TextValue "D50" '7.27'
$ws.Range("E50").Value = '  +9.70%  '

Set-TextValue "D51" '2.96'
$ws.Range("E51").Value = '  +3.89%  '

$excel.CutCopyMode = $false
